$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in this sheet store Price/Volume as text (inline strings); keep them as text
# so values like "8.650" or "0.00002106" are not mangled into numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "325.06"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-2.28%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.58"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.68%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.488"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-5.24%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08042"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-3.72%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.650"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.80%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.285"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-4.77%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.884"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-5.02%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.713"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-5.98%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9370"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.38%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1173"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-5.81%"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.72%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09939"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.55%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04255"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "8.07%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1064"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.34%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001275"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.39%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005852"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-3.93%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.40%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.537"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-5.44%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.08%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2658"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "3.35%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04247"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-3.63%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001236"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.77%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004480"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.96%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001201"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.85%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003993"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.01%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02650"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-5.77%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05495"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-3.87%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007676"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.32%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1392"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-2.48%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007399"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-17.62%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002049"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.54%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008690"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-14.62%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007109"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.57%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.01%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003528"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "8.80%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002272"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.34%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.01%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.01%"
